$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44799
$ws.Range("J2").Value = 160
$ws.Range("K2").Value = 750
$ws.Range("L2").Value = 850
$ws.Range("M2").Value = 800
$ws.Range("P2").Value = 800

# Row 3
$ws.Range("D3").Value = 44799
$ws.Range("J3").Value = 120
$ws.Range("K3").Value = 650
$ws.Range("L3").Value = 650
$ws.Range("M3").Value = 650
$ws.Range("P3").Value = 650

# Row 4
$ws.Range("D4").Value = 44797
$ws.Range("J4").Value = 240
$ws.Range("K4").Value = 750
$ws.Range("L4").Value = 850
$ws.Range("M4").Value = 800
$ws.Range("P4").Value = 800

# Row 5
$ws.Range("D5").Value = 44797
$ws.Range("K5").Value = 650
$ws.Range("L5").Value = 650
$ws.Range("M5").Value = 650
$ws.Range("P5").Value = 650

# Row 6
$ws.Range("D6").Value = 44804
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 750
$ws.Range("L6").Value = 850
$ws.Range("M6").Value = 800
$ws.Range("P6").Value = 800

# Row 7
$ws.Range("D7").Value = 44804
$ws.Range("I7").Value = "Segunda"
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 650
$ws.Range("L7").Value = 650
$ws.Range("M7").Value = 650
$ws.Range("P7").Value = 650

# Row 8
$ws.Range("D8").Value = 44791
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 240
$ws.Range("K8").Value = 750
$ws.Range("L8").Value = 800
$ws.Range("M8").Value = 775
$ws.Range("P8").Value = 775

# Row 9
$ws.Range("D9").Value = 44791
$ws.Range("I9").Value = "Segunda"
$ws.Range("J9").Value = 250
$ws.Range("K9").Value = 650
$ws.Range("L9").Value = 650
$ws.Range("M9").Value = 650
$ws.Range("P9").Value = 650

# Row 10
$ws.Range("D10").Value = 44818
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 300
$ws.Range("K10").Value = 800
$ws.Range("L10").Value = 900
$ws.Range("M10").Value = 850
$ws.Range("P10").Value = 850

# Row 11
$ws.Range("D11").Value = 44811
$ws.Range("K11").Value = 750
$ws.Range("L11").Value = 850
$ws.Range("M11").Value = 800
$ws.Range("P11").Value = 800

# Row 12
$ws.Range("D12").Value = 44859
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 300
$ws.Range("K12").Value = 700
$ws.Range("L12").Value = 800
$ws.Range("M12").Value = 750
$ws.Range("P12").Value = 750

# Row 13
$ws.Range("D13").Value = 44859
$ws.Range("I13").Value = "Segunda"
$ws.Range("J13").Value = 200
$ws.Range("L13").Value = 600
$ws.Range("M13").Value = 600
$ws.Range("P13").Value = 600

# Row 14
$ws.Range("D14").Value = 44624
$ws.Range("J14").Value = 120
$ws.Range("K14").Value = 650
$ws.Range("L14").Value = 700
$ws.Range("M14").Value = 675
$ws.Range("P14").Value = 675

# Row 15
$ws.Range("D15").Value = 44831
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 300
$ws.Range("K15").Value = 700
$ws.Range("L15").Value = 800
$ws.Range("M15").Value = 750
$ws.Range("P15").Value = 750

# Row 16
$ws.Range("D16").Value = 44831
$ws.Range("I16").Value = "Segunda"
$ws.Range("J16").Value = 200
$ws.Range("K16").Value = 600
$ws.Range("L16").Value = 600
$ws.Range("M16").Value = 600
$ws.Range("P16").Value = 600

# Row 17
$ws.Range("D17").Value = 44837
$ws.Range("K17").Value = 700
$ws.Range("L17").Value = 800
$ws.Range("M17").Value = 750
$ws.Range("P17").Value = 750

# Row 18
$ws.Range("D18").Value = 44837
$ws.Range("J18").Value = 150
$ws.Range("K18").Value = 600
$ws.Range("L18").Value = 600
$ws.Range("M18").Value = 600
$ws.Range("P18").Value = 600

# Row 19
$ws.Range("D19").Value = 44764
$ws.Range("J19").Value = 200

# Row 20
$ws.Range("D20").Value = 44764
$ws.Range("J20").Value = 150

# Row 21
$ws.Range("D21").Value = 44883
$ws.Range("J21").Value = 300

# Row 22
$ws.Range("D22").Value = 44883
$ws.Range("J22").Value = 200

# Row 23
$ws.Range("D23").Value = 44882
$ws.Range("J23").Value = 400
$ws.Range("K23").Value = 700
$ws.Range("M23").Value = 750
$ws.Range("P23").Value = 750

# Row 24
$ws.Range("D24").Value = 44882
$ws.Range("J24").Value = 300
$ws.Range("K24").Value = 600
$ws.Range("L24").Value = 600
$ws.Range("M24").Value = 600
$ws.Range("P24").Value = 600

# Row 25
$ws.Range("D25").Value = 44608
$ws.Range("J25").Value = 120
$ws.Range("K25").Value = 600
$ws.Range("L25").Value = 650
$ws.Range("M25").Value = 625
$ws.Range("P25").Value = 625

# Row 26
$ws.Range("D26").Value = 44839
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 240
$ws.Range("K26").Value = 700
$ws.Range("L26").Value = 800
$ws.Range("M26").Value = 750
$ws.Range("P26").Value = 750

# Row 27
$ws.Range("D27").Value = 44839
$ws.Range("I27").Value = "Segunda"
$ws.Range("K27").Value = 600
$ws.Range("L27").Value = 600
$ws.Range("M27").Value = 600
$ws.Range("P27").Value = 600

# Row 28
$ws.Range("D28").Value = 44761
$ws.Range("J28").Value = 200
$ws.Range("K28").Value = 700
$ws.Range("L28").Value = 800
$ws.Range("M28").Value = 750
$ws.Range("P28").Value = 750

# Row 29
$ws.Range("D29").Value = 44761
$ws.Range("I29").Value = "Segunda"
$ws.Range("J29").Value = 150
$ws.Range("L29").Value = 600
$ws.Range("M29").Value = 600
$ws.Range("P29").Value = 600

# Row 30
$ws.Range("D30").Value = 44868
$ws.Range("L30").Value = 800
$ws.Range("M30").Value = 750
$ws.Range("P30").Value = 750

# Row 31
$ws.Range("D31").Value = 44754
$ws.Range("L31").Value = 750
$ws.Range("M31").Value = 725
$ws.Range("P31").Value = 725

# Row 32
$ws.Range("D32").Value = 44610
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 100
$ws.Range("L32").Value = 650
$ws.Range("M32").Value = 625
$ws.Range("P32").Value = 625
